$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing header row 3: D3/E3 should be bold like B3/C3 (style s=1) ---
$ws.Range("D3:E3").Font.Bold = $true

# --- Remove the stray leftover cell at C17 (old italic-style placeholder) ---
$ws.Rows("17").Delete()

# =========================================================
# Table: Best r2 scores Stockholm (rows 13-21)
# =========================================================
$ws.Range("B13").Value = "Best r2 scores Stockholm"
$ws.Range("B13:C13").Font.Bold = $true
$ws.Range("B13:C13").HorizontalAlignment = -4108
$ws.Range("B13:C13").Merge() | Out-Null

$ws.Range("B14").Value = "Model"
$ws.Range("C14").Value = "Eval method"
$ws.Range("D14").Value = "Train r2"
$ws.Range("E14").Value = "Test r2"
$ws.Range("B14:E14").Font.Bold = $true

$ws.Range("B15").Value = "Linear Regression"
$ws.Range("C15").Value = "LOOCV"
$ws.Range("E15").Value = 0.574

$ws.Range("B16").Value = "Linear Regression"
$ws.Range("C16").Value = "Average score"
$ws.Range("D16").Value = 0.604
$ws.Range("D5").Copy()
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Value = 0.5559

$ws.Range("B17").Value = "NuSVR"
$ws.Range("C17").Value = "Average score"
$ws.Range("D17").Value = 0.7208
$ws.Range("E17").Value = 0.5857

$ws.Range("B18").Value = "Gradient Boosting"
$ws.Range("C18").Value = "Average score"
$ws.Range("D18").Value = 0.9434
$ws.Range("E18").Value = 0.606

$ws.Range("B19").Value = "Random Forest"
$ws.Range("C19").Value = "Average score"
$ws.Range("D19").Value = 0.8764
$ws.Range("E19").Value = 0.6025

$ws.Range("B20").Value = "XGBoost"
$ws.Range("C20").Value = "Average score"
$ws.Range("D20").Value = 0.9811
$ws.Range("E20").Value = 0.6175

$ws.Range("B21").Value = "Neural Network"
$ws.Range("C21").Value = "Average score"
$ws.Range("D21").Value = 0.6335
$ws.Range("E21").Value = 0.5905

# =========================================================
# Table: London (rows 24-27)
# =========================================================
$ws.Range("B24").Value = "London"

$ws.Range("B25").Value = "Model"
$ws.Range("C25").Value = "Eval method"
$ws.Range("D25").Value = "Train r2"
$ws.Range("E25").Value = "Test r2"
$ws.Range("B25:E25").Font.Bold = $true

$ws.Range("B26").Value = "Gradient Boosting"
$ws.Range("C26").Value = "Average score"
$ws.Range("D26").Value = 0.8676
$ws.Range("E26").Value = 0.6596

$ws.Range("B27").Value = "XGBoost"
$ws.Range("C27").Value = "Average score"
$ws.Range("D27").Value = 0.797
$ws.Range("E27").Value = 0.6455

# =========================================================
# Table: Paris (rows 29-32)
# =========================================================
$ws.Range("B29").Value = "Paris"

$ws.Range("B30").Value = "Model"
$ws.Range("C30").Value = "Eval method"
$ws.Range("D30").Value = "Train r2"
$ws.Range("E30").Value = "Test r2"
$ws.Range("B30:E30").Font.Bold = $true

$ws.Range("B31").Value = "Gradient Boosting"
$ws.Range("C31").Value = "Average score"
$ws.Range("D31").Value = 0.8536
$ws.Range("E31").Value = 0.7409

$ws.Range("B32").Value = "XGBoost"
$ws.Range("C32").Value = "Average score"
$ws.Range("D32").Value = 0.807
$ws.Range("E32").Value = 0.7052

# --- Selection matches final diff state ---
$ws.Range("E33").Select() | Out-Null

Write-Host "edit complete"
